$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 426, shifting the existing rows (426..442)
# down to (428..444). Excel's Range.Insert() shifts cells down and keeps the
# row's number-format (style) carried on column D.
$ws.Range("A426:R427").Insert()

# Fill the two newly-inserted rows with their final values.
$ws.Range("A426").Value = 7
$ws.Range("B426").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C426").Value = "Ñuble"
$ws.Range("D426").Value = 45075
$ws.Range("E426").Value = 16
$ws.Range("F426").Value = 100114013
$ws.Range("G426").Value = "Zanahoria"
$ws.Range("H426").Value = "Sin especificar"
$ws.Range("I426").Value = "Primera"
$ws.Range("J426").Value = 80
$ws.Range("K426").Value = 7000
$ws.Range("L426").Value = 7000
$ws.Range("M426").Value = 7000
$ws.Range("N426").Value = "$/saco 20 kilos"
$ws.Range("O426").Value = "Provincia de Diguillín"
$ws.Range("P426").Value = 350
$ws.Range("Q426").Value = 20
$ws.Range("R426").Value = "Hortaliza"

$ws.Range("A427").Value = 7
$ws.Range("B427").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C427").Value = "Ñuble"
$ws.Range("D427").Value = 45075
$ws.Range("E427").Value = 16
$ws.Range("F427").Value = 100114013
$ws.Range("G427").Value = "Zanahoria"
$ws.Range("H427").Value = "Sin especificar"
$ws.Range("I427").Value = "Segunda"
$ws.Range("J427").Value = 60
$ws.Range("K427").Value = 6000
$ws.Range("L427").Value = 6000
$ws.Range("M427").Value = 6000
$ws.Range("N427").Value = "$/saco 20 kilos"
$ws.Range("O427").Value = "Provincia de Diguillín"
$ws.Range("P427").Value = 300
$ws.Range("Q427").Value = 20
$ws.Range("R427").Value = "Hortaliza"
